# Completed Bro Code favicon tutorial: replace the leftover <video> attribute
# list (src/width/controls/autoplay/muted/loop) with the new "Favicons"
# section (Favicons header, <link> bullet, rel bullet) -- including the
# _GoBack bookmark that sat inside the old "loop" paragraph, now relocated
# inside the new "rel" paragraph.

$d = $word.ActiveDocument

# Locate the contiguous run of paragraphs that make up the old "video
# attributes" block by matching on their text, rather than hard-coding
# indices, so the script is resilient to minor paragraph-count drift.
$startIndex = -1
$endIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "src = defines the URL of source video*") {
        $startIndex = $i
    }
    if ($t -like "loop = loops video*") {
        $endIndex = $i
    }
}

if ($startIndex -eq -1 -or $endIndex -eq -1 -or $endIndex -lt $startIndex) {
    throw "Could not locate the video-attributes paragraph block to replace"
}

$startPara = $d.Paragraphs.Item($startIndex)
$endPara = $d.Paragraphs.Item($endIndex)

# Remove the whole old block (src / width / controls / autoplay / muted / loop).
$oldRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$oldRange.Delete()

# Insert a fresh paragraph right after the preceding paragraph (the
# "<video></video> = ..." line) to host the new content, then replace it
# (via InsertXML) with the three new paragraphs in one shot.
$anchorPara = $d.Paragraphs.Item($startIndex - 1)
$anchorRange = $anchorPara.Range
$anchorRange.Collapse(0)
$anchorRange.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($startIndex)
$newRange = $newPara.Range
$newRange.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:i/><w:u w:val="single"/></w:rPr><w:t>Favicons</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>&lt;link&gt;</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">= </w:t></w:r><w:r><w:t>defines relationship between current document and external resource</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>rel</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = </w:t></w:r><w:r><w:t>defines relationship between current docume</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>nt and external resource</w:t></w:r></w:p>
'@)

Write-Output "Replaced paragraphs $startIndex..$endIndex with Favicons/link/rel block; new paragraph count: $($d.Paragraphs.Count)"
